$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new data row at row 117, pushing the existing rows 117:222
# down to 118:223 (matches the dimension change A1:R222 -> A1:R223).
$ws.Rows.Item(117).Insert()

# Populate the newly inserted row 117 with the new record.
$ws.Range("A117").Value = 7
$ws.Range("B117").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C117").Value = "Ñuble"
$ws.Range("D117").Value = 44587
$ws.Range("E117").Value = 16
$ws.Range("F117").Value = 100112008
$ws.Range("G117").Value = "Coliflor"
$ws.Range("H117").Value = "Sin especificar"
$ws.Range("I117").Value = "Primera"
$ws.Range("J117").Value = 200
$ws.Range("K117").Value = 800
$ws.Range("L117").Value = 850
$ws.Range("M117").Value = 825
$ws.Range("N117").Value = "$/unidad"
$ws.Range("O117").Value = "Provincia de Diguillín"
$ws.Range("P117").Value = 825
$ws.Range("Q117").Value = 1
$ws.Range("R117").Value = "Hortaliza"
